{"js": "// Zulu translation corrections (5 textual fixes) for the ParentText\n// 5-Day UX RCT Quantitative Participant Information Sheets and Consent\n// Form document.\nconst body = context.document.body;\n\nconst replacements = [\n  {\n    find: \"sicela uthumele i-email ithimbeni locwaningo ku-\",\n    replace: \"sicela uthumele i-email ithimba locwaningo ku-\",\n  },\n  {\n    find: \"kuqinisekise ukuthi uzizwa ukhululekile uma uphendula imibuzo.\",\n    replace: \"kuqinisekise ukuthi uzizwe ukhululekile uma uphendula imibuzo.\",\n  },\n  {\n    find: \"kanti iMenenja yocwaningo nguZamakhanya Makhanya (University of Cape Town).\",\n    replace: \"kanye neMenenja yocwaningo uZamakhanya Makhanya (University of Cape Town).\",\n  },\n  {\n    find: \"Uma unemibuzo noma okukukhathazayo mayelana namalungelo akho njengomhlanganyeli wocwaningo, ungathintana nethimba locwaningo ku-\",\n    replace: \"Uma unemibuzo noma kukhona okukukhathazayo mayelana namalungelo akho njengomhlanganyeli wocwaningo, ungathintana nethimba locwaningo ku-\",\n  },\n  {\n    find: \"Uma ufunde futhi waqonda idokhumenti engenhla, vuma kulemilayezo\",\n    replace: \"Uma ufunde futhi waqonda incwadi engenhla, vuma kulemilayezo\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Zulu translation corrections (5 textual fixes) for the ParentText\n# 5-Day UX RCT Quantitative Participant Information Sheets and Consent\n# Form document.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"sicela uthumele i-email ithimbeni locwaningo ku-\"; Replace = \"sicela uthumele i-email ithimba locwaningo ku-\" },\n    @{ Find = \"kuqinisekise ukuthi uzizwa ukhululekile uma uphendula imibuzo.\"; Replace = \"kuqinisekise ukuthi uzizwe ukhululekile uma uphendula imibuzo.\" },\n    @{ Find = \"kanti iMenenja yocwaningo nguZamakhanya Makhanya (University of Cape Town).\"; Replace = \"kanye neMenenja yocwaningo uZamakhanya Makhanya (University of Cape Town).\" },\n    @{ Find = \"Uma unemibuzo noma okukukhathazayo mayelana namalungelo akho njengomhlanganyeli wocwaningo, ungathintana nethimba locwaningo ku-\"; Replace = \"Uma unemibuzo noma kukhona okukukhathazayo mayelana namalungelo akho njengomhlanganyeli wocwaningo, ungathintana nethimba locwaningo ku-\" },\n    @{ Find = \"Uma ufunde futhi waqonda idokhumenti engenhla, vuma kulemilayezo\"; Replace = \"Uma ufunde futhi waqonda incwadi engenhla, vuma kulemilayezo\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $found = $find.Execute($r.Find, $false, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n    if (-not $found) {\n        throw \"Text not found: $($r.Find)\"\n    }\n}\n"}
